# "101 Pass for KLIC"
# - Update the "TestCaseNumber=102" value to "TestCaseNumber=101" on the
#   "Test Cases" sheet (cell D2).
# - Move the window's first-visible-tab / selection so that the "Test
#   Cases" sheet is scrolled into view and the active selection on that
#   sheet moves from C3:C6 to D3.

$wb = $excel.ActiveWorkbook

$testCasesSheet = $wb.Worksheets.Item("Test Cases")

# Update the TestCaseNumber value in D2 (102 -> 101)
$testCasesSheet.Range("D2").Value = "TestCaseNumber=101"

# Scroll the workbook tabs so the "Test Cases" sheet (the active one) is
# the first visible tab in the window.
$testCasesSheet.Activate()
$wb.Windows.Item(1).ScrollWorkbookTabs(1)

# Update the selection on the "Test Cases" sheet from C3:C6 to D3
$testCasesSheet.Range("D3").Select()
